# edits for petB primer paper
# Updates the colour palette on the "subclade" sheet: several subclade /
# colour_name pairs are recoloured, and a new "Others" / "#4D4D4D" row is
# appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("subclade")

# NOTE: the order these new colour strings are first written in matters for
# shared-string layout, so touch the brand-new values (row 36, then the
# genuinely-new hex codes) before re-using colours that already exist
# elsewhere in the workbook.
$ws.Range("B36").Value = "#4D4D4D"
$ws.Range("A36").Value = "Others"

$ws.Range("B23").Value = "#51A3CC"
$ws.Range("B8").Value  = "#CC5151"
$ws.Range("B4").Value  = "#A3CC51"
$ws.Range("B5").Value  = "#A3CC51"
$ws.Range("B6").Value  = "#6B990F"
$ws.Range("B22").Value = "#7EC3E5"

# Re-used existing colours.
$ws.Range("B7").Value  = "#E57E7E"
$ws.Range("B25").Value = "#FF3300"
$ws.Range("B28").Value = "#E57E7E"
$ws.Range("B29").Value = "#FFB2B2"
$ws.Range("B32").Value = "#BFB2FF"

# Matches the saved selection/scroll state recorded for this sheet after the
# edit (topLeftCell scroll anchor cleared, selection moved to E5).
$ws.Range("E5").Select() | Out-Null
